# Rename header cells on the existing sheets.
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "PO Forecast"

# Header row.
$ws.Cells.Item(1,1).Value = "ds"
$ws.Cells.Item(1,2).Value = "PO_Forecast"
$ws.Cells.Item(1,3).Value = "yhat_lower"
$ws.Cells.Item(1,4).Value = "yhat_upper"

# Match the bold/bordered header style used on the other sheets.
$wsWeekly.Range("B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Match the date-formatted style used for the "ds" / date column.
$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A24").PasteSpecial(-4122)

# Data rows.

$ws.Cells.Item(2,1).Value = 45333.99999999999
$ws.Cells.Item(2,2).Value = 18
$ws.Cells.Item(2,3).Value = -4.263891294227465
$ws.Cells.Item(2,4).Value = 38.98598527288112
$ws.Cells.Item(3,1).Value = 45382.99999999999
$ws.Cells.Item(3,2).Value = 17
$ws.Cells.Item(3,3).Value = -3.934442629014511
$ws.Cells.Item(3,4).Value = 38.47876469668536
$ws.Cells.Item(4,1).Value = 45424.99999999999
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = -5.447905823390353
$ws.Cells.Item(4,4).Value = 37.67288196493281
$ws.Cells.Item(5,1).Value = 45438.99999999999
$ws.Cells.Item(5,2).Value = 17
$ws.Cells.Item(5,3).Value = -5.692827671838336
$ws.Cells.Item(5,4).Value = 39.05388467955598
$ws.Cells.Item(6,1).Value = 45445.99999999999
$ws.Cells.Item(6,2).Value = 17
$ws.Cells.Item(6,3).Value = -5.674146616734724
$ws.Cells.Item(6,4).Value = 38.05097022613671
$ws.Cells.Item(7,1).Value = 45452.99999999999
$ws.Cells.Item(7,2).Value = 17
$ws.Cells.Item(7,3).Value = -4.760008887520162
$ws.Cells.Item(7,4).Value = 38.70122591154396
$ws.Cells.Item(8,1).Value = 45466.99999999999
$ws.Cells.Item(8,2).Value = 16
$ws.Cells.Item(8,3).Value = -3.663954392779537
$ws.Cells.Item(8,4).Value = 39.74932804002853
$ws.Cells.Item(9,1).Value = 45473.99999999999
$ws.Cells.Item(9,2).Value = 16
$ws.Cells.Item(9,3).Value = -4.269622396434148
$ws.Cells.Item(9,4).Value = 37.67966356708929
$ws.Cells.Item(10,1).Value = 45480.99999999999
$ws.Cells.Item(10,2).Value = 16
$ws.Cells.Item(10,3).Value = -5.607961722057667
$ws.Cells.Item(10,4).Value = 39.12374312277407
$ws.Cells.Item(11,1).Value = 45494.99999999999
$ws.Cells.Item(11,2).Value = 16
$ws.Cells.Item(11,3).Value = -6.801232713547844
$ws.Cells.Item(11,4).Value = 38.06621166907205
$ws.Cells.Item(12,1).Value = 45515.99999999999
$ws.Cells.Item(12,2).Value = 16
$ws.Cells.Item(12,3).Value = -4.658908340081516
$ws.Cells.Item(12,4).Value = 36.27615239988178
$ws.Cells.Item(13,1).Value = 45536.99999999999
$ws.Cells.Item(13,2).Value = 16
$ws.Cells.Item(13,3).Value = -7.743983178212615
$ws.Cells.Item(13,4).Value = 37.72022345083754
$ws.Cells.Item(14,1).Value = 45543.99999999999
$ws.Cells.Item(14,2).Value = 16
$ws.Cells.Item(14,3).Value = -6.624281563708246
$ws.Cells.Item(14,4).Value = 37.13630588437487
$ws.Cells.Item(15,1).Value = 45550.99999999999
$ws.Cells.Item(15,2).Value = 16
$ws.Cells.Item(15,3).Value = -5.202400788525252
$ws.Cells.Item(15,4).Value = 38.81713401650698
$ws.Cells.Item(16,1).Value = 45557.99999999999
$ws.Cells.Item(16,2).Value = 16
$ws.Cells.Item(16,3).Value = -5.243257878796212
$ws.Cells.Item(16,4).Value = 37.17369731007807
$ws.Cells.Item(17,1).Value = 45564.99999999999
$ws.Cells.Item(17,2).Value = 16
$ws.Cells.Item(17,3).Value = -4.145666156503219
$ws.Cells.Item(17,4).Value = 36.77852673418747
$ws.Cells.Item(18,1).Value = 45571.99999999999
$ws.Cells.Item(18,2).Value = 15
$ws.Cells.Item(18,3).Value = -6.869744845221506
$ws.Cells.Item(18,4).Value = 36.22211643803022
$ws.Cells.Item(19,1).Value = 45578.99999999999
$ws.Cells.Item(19,2).Value = 15
$ws.Cells.Item(19,3).Value = -6.530477528606498
$ws.Cells.Item(19,4).Value = 35.75980092730923
$ws.Cells.Item(20,1).Value = 45585.99999999999
$ws.Cells.Item(20,2).Value = 15
$ws.Cells.Item(20,3).Value = -6.489552176113326
$ws.Cells.Item(20,4).Value = 36.70417694295122
$ws.Cells.Item(21,1).Value = 45592.99999999999
$ws.Cells.Item(21,2).Value = 15
$ws.Cells.Item(21,3).Value = -5.56456181009367
$ws.Cells.Item(21,4).Value = 36.98911847109493
$ws.Cells.Item(22,1).Value = 45599.99999999999
$ws.Cells.Item(22,2).Value = 15
$ws.Cells.Item(22,3).Value = -6.136198068616865
$ws.Cells.Item(22,4).Value = 37.541776299582
$ws.Cells.Item(23,1).Value = 45606.99999999999
$ws.Cells.Item(23,2).Value = 15
$ws.Cells.Item(23,3).Value = -6.879432145215072
$ws.Cells.Item(23,4).Value = 34.23146102499184
$ws.Cells.Item(24,1).Value = 45613.99999999999
$ws.Cells.Item(24,2).Value = 15
$ws.Cells.Item(24,3).Value = -7.930183164975018
$ws.Cells.Item(24,4).Value = 35.49243074629586

$excel.CutCopyMode = $false
